$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C (old C..G shift to D..H).
$ws.Columns("C").Insert() | Out-Null

# New column C gets the same width as column B (~18.29 chars wide), no bestFit.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Give the new column's data cells (rows 2-20) the bordered / centered look
# used throughout the rest of the table.
$dataRange = $ws.Range("C2:C20")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.Borders.LineStyle = 1

# Populate the "module exists" marker ("v") for the relevant rows.
$ws.Range("C3").Value = "v"
$ws.Range("C7").Value = "v"
$ws.Range("C11").Value = "v"
$ws.Range("C13").Value = "v"
$ws.Range("C14").Value = "v"
$ws.Range("C16").Value = "v"

# Move the active selection to C14, matching the saved view state.
$ws.Range("C14").Select() | Out-Null
